# Slide 10 ("Thank You"): grow the "Content Placeholder 2" shape upward and
# add an empty-line gap plus a new Github credit line under "ExModelo Team".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(10)
$sh = $s.Shapes.Item(2)

# Re-position / re-size the placeholder shape (it's taller now to fit the
# extra lines of text).
$sh.Top    = 264.3301
$sh.Height = 175.7675590551181

# Append two blank paragraphs and a new "Github: pioucyril/exmodeloGroup4"
# paragraph after the existing "Acknowledgement: ExModelo Team" text, while
# preserving the pre-existing runs/formatting. Each InsertAfter call targets
# the live end-of-text range so the new runs land in the right place and
# stay split the way the source file splits them.
$tr = $sh.TextFrame.TextRange
[void]$tr.InsertAfter("`r`r`r")

$tr = $sh.TextFrame.TextRange
[void]$tr.InsertAfter("Github")

$tr = $sh.TextFrame.TextRange
[void]$tr.InsertAfter(": ")

$tr = $sh.TextFrame.TextRange
[void]$tr.InsertAfter("pioucyril")

$tr = $sh.TextFrame.TextRange
[void]$tr.InsertAfter("/exmodeloGroup4")
